# Adding new code 05.04.2020
$wb = $excel.ActiveWorkbook

# --- "files" worksheet: insert a new "Execution" column before "Result" ---
$filesWs = $wb.Worksheets.Item("files")

# Insert the new column at H (inside the existing F:H 28.109375-wide block) so
# the engine extends that column-width span to F:I, exactly like Excel does
# when a new column lands in the middle of a uniform-width run. This shifts
# the old H (Target_key) into I and the old I (Result) into J.
$filesWs.Columns.Item(8).Insert()

# Put Target_key/BusinessID back where they belong (column H) and reset the
# formatting the insert left behind (it inherited the neighbouring wrap style).
$filesWs.Cells.Item(1, 8).Value = "Target_key"
$filesWs.Cells.Item(2, 8).Value = "BusinessID"
$filesWs.Cells.Item(2, 8).Style = "Normal"

# The brand-new "Execution" column now correctly sits at I.
$filesWs.Cells.Item(1, 9).Value = "Execution"
$filesWs.Cells.Item(2, 9).Value = "Y"
$filesWs.Cells.Item(2, 9).WrapText = $true

# New data row (row 3), mirroring row 2 with an incremented TC_No.
$filesWs.Cells.Item(3, 1).Value = 2
$filesWs.Cells.Item(3, 2).Value = "File_testing"
$filesWs.Cells.Item(3, 3).Value = "D://Database_automation//Database_automation//files//client_data_src.csv"
$filesWs.Cells.Item(3, 4).Value = ","
$filesWs.Cells.Item(3, 5).Value = "BusinessID"
$filesWs.Cells.Item(3, 6).Value = "D://Database_automation//Database_automation//files//client_data_target.csv"
$filesWs.Cells.Item(3, 7).Value = ","
$filesWs.Cells.Item(3, 8).Value = "BusinessID"
$filesWs.Cells.Item(3, 9).Value = "Y"

$filesWs.Range("C3:D3").WrapText = $true
$filesWs.Range("F3:G3").WrapText = $true
$filesWs.Cells.Item(3, 9).WrapText = $true
$filesWs.Rows.Item(3).RowHeight = $filesWs.Rows.Item(2).RowHeight

$filesWs.Range("C10").Select()

# --- "database" worksheet: just a cursor move, no data changes ---
$dbWs = $wb.Worksheets.Item("database")
$dbWs.Range("E2").Select()

$filesWs.Activate()
